# The sighting rows (3-17) on the sheet got reshuffled into a different row
# order (same 15 observations, new row positions). Capture the current
# per-row values for the columns that move with the record (species/taxon
# fields, coordinates, and the collector-order text), then write them back
# out in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 17
$cols = @("A","B","D","E","F","G","H","Q","R","AX")

# New row order: $newOrder[$targetRow] = $sourceRow that its data comes from.
$newOrder = @{
    3  = 13
    4  = 11
    5  = 3
    6  = 5
    7  = 17
    8  = 6
    9  = 8
    10 = 12
    11 = 14
    12 = 15
    13 = 7
    14 = 16
    15 = 4
    16 = 9
    17 = 10
}

# Snapshot every source row's values before any writes happen, so a row we
# already overwrote never gets read back out as a source for another row.
$snapshot = @{}
foreach ($r in $firstRow..$lastRow) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($targetRow in $firstRow..$lastRow) {
    $sourceRow = $newOrder[$targetRow]
    $rowVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $rowVals[$c]
    }
}
